# Daily attendance processing - 2025-10-12 07:40:33
# Reorder the "Recorded By" audit-trail text in column G for specific rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose "Recorded By" value is exactly "dnasr281@gmail.com, System"
# and needs to become "System, dnasr281@gmail.com"
$dnasrRows = @(3,6,10,11,12,13,14,15,17,30,33,37,38,39,40,41,42,44,57,60,64,65,66,67,68,69,71,86,87,88,89,93,95,96,112,113,114,115,119,121,122,138,139,140,141,145,147,148)

foreach ($r in $dnasrRows) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq "dnasr281@gmail.com, System") {
        $cell.Value = "System, dnasr281@gmail.com"
    }
}

# Rows whose "Recorded By" value is exactly "backup@backdoor.com, system, System"
# and needs to become "backup@backdoor.com, System, system"
$backupRows = @(2,29,56)

foreach ($r in $backupRows) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq "backup@backdoor.com, system, System") {
        $cell.Value = "backup@backdoor.com, System, system"
    }
}
